$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "643×6=3858"
$t.Cell(1,2).Range.Text  = "425×7=2975"
$t.Cell(1,3).Range.Text  = "116×2=232"
$t.Cell(1,4).Range.Text  = "779×2=1558"
$t.Cell(1,5).Range.Text  = "937×7=6559"

$t.Cell(5,1).Range.Text  = "717×4=2868"
$t.Cell(5,2).Range.Text  = "570×2=1140"
$t.Cell(5,3).Range.Text  = "417×6=2502"
$t.Cell(5,4).Range.Text  = "136×9=1224"
$t.Cell(5,5).Range.Text  = "360×8=2880"

$t.Cell(10,1).Range.Text = "659×4=2636"
$t.Cell(10,2).Range.Text = "930×2=1860"
$t.Cell(10,3).Range.Text = "199×7=1393"
$t.Cell(10,4).Range.Text = "545×6=3270"
$t.Cell(10,5).Range.Text = "658×9=5922"

$t.Cell(15,1).Range.Text = "929×4=3716"
$t.Cell(15,2).Range.Text = "435×8=3480"
$t.Cell(15,3).Range.Text = "214×8=1712"
$t.Cell(15,4).Range.Text = "867×6=5202"
$t.Cell(15,5).Range.Text = "339×9=3051"

$t.Cell(20,1).Range.Text = "428×6=2568"
$t.Cell(20,2).Range.Text = "518×7=3626"
$t.Cell(20,3).Range.Text = "908×2=1816"
$t.Cell(20,4).Range.Text = "846×6=5076"
$t.Cell(20,5).Range.Text = "458×9=4122"
